# "Fruta / hortaliza, semanal" - weekly Cebollín price records for
# Terminal Hortofrutícola Agro Chillán: re-date/re-value the existing
# daily rows (44-62) to their weekly-aggregated figures and append two
# new rows (63-64) for the records that were pushed out of the range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45
$ws.Range("D45").Value = 44917
$ws.Range("J45").Value = 400
$ws.Range("K45").Value = 700
$ws.Range("L45").Value = 700
$ws.Range("M45").Value = 700
$ws.Range("N45").Value = "`$/paquete 6 unidades"
$ws.Range("P45").Value = 117
$ws.Range("Q45").Value = 6

# Row 46
$ws.Range("D46").Value = 44917
$ws.Range("I46").Value = "Segunda"
$ws.Range("J46").Value = 300
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 600
$ws.Range("M46").Value = 600
$ws.Range("N46").Value = "`$/paquete 6 unidades"
$ws.Range("P46").Value = 100
$ws.Range("Q46").Value = 6

# Row 47
$ws.Range("D47").Value = 44775
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("M47").Value = 8000
$ws.Range("N47").Value = "`$/docena de atados"
$ws.Range("P47").Value = 2667
$ws.Range("Q47").Value = 3

# Row 48
$ws.Range("D48").Value = 44782
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 100
$ws.Range("K48").Value = 8000
$ws.Range("L48").Value = 8000
$ws.Range("M48").Value = 8000
$ws.Range("N48").Value = "`$/docena de atados"
$ws.Range("P48").Value = 2667
$ws.Range("Q48").Value = 3

# Row 49
$ws.Range("D49").Value = 44894
$ws.Range("J49").Value = 400
$ws.Range("K49").Value = 600
$ws.Range("L49").Value = 700
$ws.Range("M49").Value = 650
$ws.Range("N49").Value = "`$/paquete 6 unidades"
$ws.Range("P49").Value = 108
$ws.Range("Q49").Value = 6

# Row 50
$ws.Range("D50").Value = 44894
$ws.Range("I50").Value = "Segunda"
$ws.Range("J50").Value = 300
$ws.Range("K50").Value = 500
$ws.Range("L50").Value = 500
$ws.Range("M50").Value = 500
$ws.Range("N50").Value = "`$/paquete 6 unidades"
$ws.Range("P50").Value = 83
$ws.Range("Q50").Value = 6

# Row 51
$ws.Range("D51").Value = 44769
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = 8000
$ws.Range("N51").Value = "`$/docena de atados"
$ws.Range("P51").Value = 2667
$ws.Range("Q51").Value = 3

# Row 52
$ws.Range("D52").Value = 44791
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 120
$ws.Range("K52").Value = 8000
$ws.Range("L52").Value = 8500
$ws.Range("M52").Value = 8250
$ws.Range("N52").Value = "`$/docena de atados"
$ws.Range("P52").Value = 2750
$ws.Range("Q52").Value = 3

# Row 53
$ws.Range("J53").Value = 300

# Row 55
$ws.Range("D55").Value = 44895
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 600
$ws.Range("L55").Value = 700
$ws.Range("M55").Value = 650
$ws.Range("N55").Value = "`$/paquete 6 unidades"
$ws.Range("P55").Value = 108
$ws.Range("Q55").Value = 6

# Row 56
$ws.Range("D56").Value = 44895
$ws.Range("I56").Value = "Segunda"
$ws.Range("J56").Value = 300
$ws.Range("K56").Value = 500
$ws.Range("L56").Value = 500
$ws.Range("M56").Value = 500
$ws.Range("N56").Value = "`$/paquete 6 unidades"
$ws.Range("P56").Value = 83
$ws.Range("Q56").Value = 6

# Row 57
$ws.Range("D57").Value = 44847
$ws.Range("K57").Value = 7500
$ws.Range("L57").Value = 8000
$ws.Range("M57").Value = 7750
$ws.Range("P57").Value = 2583

# Row 58
$ws.Range("D58").Value = 44839
$ws.Range("J58").Value = 120
$ws.Range("K58").Value = 7500
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = 7750
$ws.Range("N58").Value = "`$/docena de atados"
$ws.Range("P58").Value = 2583
$ws.Range("Q58").Value = 3

# Row 59
$ws.Range("D59").Value = 44818
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 120
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 8500
$ws.Range("M59").Value = 8250
$ws.Range("N59").Value = "`$/docena de atados"
$ws.Range("P59").Value = 2750
$ws.Range("Q59").Value = 3

# Row 60
$ws.Range("D60").Value = 44883

# Row 61
$ws.Range("D61").Value = 44883

# Row 62
$ws.Range("D62").Value = 44879
$ws.Range("J62").Value = 400
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 700
$ws.Range("M62").Value = 650
$ws.Range("N62").Value = "`$/paquete 6 unidades"
$ws.Range("O62").Value = "Provincia de Diguillín"
$ws.Range("P62").Value = 108
$ws.Range("Q62").Value = 6

# Row 63 (new)
$ws.Range("A63").Value = 7
$ws.Range("B63").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C63").Value = "Ñuble"
$ws.Range("D63").Value = 44879
$ws.Range("D63").NumberFormat = $ws.Range("D62").NumberFormat
$ws.Range("E63").Value = 16
$ws.Range("F63").Value = 100112037
$ws.Range("G63").Value = "Cebollín"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Segunda"
$ws.Range("J63").Value = 300
$ws.Range("K63").Value = 500
$ws.Range("L63").Value = 500
$ws.Range("M63").Value = 500
$ws.Range("N63").Value = "`$/paquete 6 unidades"
$ws.Range("O63").Value = "Provincia de Diguillín"
$ws.Range("P63").Value = 83
$ws.Range("Q63").Value = 6
$ws.Range("R63").Value = "Hortaliza"

# Row 64 (new)
$ws.Range("A64").Value = 7
$ws.Range("B64").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C64").Value = "Ñuble"
$ws.Range("D64").Value = 44762
$ws.Range("D64").NumberFormat = $ws.Range("D62").NumberFormat
$ws.Range("E64").Value = 16
$ws.Range("F64").Value = 100112037
$ws.Range("G64").Value = "Cebollín"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 60
$ws.Range("K64").Value = 8000
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = 8000
$ws.Range("N64").Value = "`$/docena de atados"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("P64").Value = 2667
$ws.Range("Q64").Value = 3
$ws.Range("R64").Value = "Hortaliza"
